# Adds the "Week6" journal entry block before the trailing placeholder
# paragraph that currently holds the "_GoBack" bookmark (mirrors the
# structure used for every previous "WeekN" section in this journal).

$d = $word.ActiveDocument

# Locate the (normally invisible) "_GoBack" bookmark left over from the
# last edit in Word - its paragraph is the empty placeholder at the very
# end of the document body, right after the "Week5" section.
$bm = $d.Bookmarks("_GoBack")
$target = $bm.Range.Paragraphs(1).Range
$target.Collapse(1)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Week6</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Previous period of work completed</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>Andrew Campbell, 4 hours.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>David Doherty, 4 hours.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>John Lee Heaney</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>,  8</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> hours.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>Matthew Mc Clean, 4 hours.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Ademolu</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> - 1hr</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:color w:val="8496B0" w:themeColor="text2" w:themeTint="99"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:color w:val="8496B0" w:themeColor="text2" w:themeTint="99"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Evaluation of outstanding items of previous period of work or additional work completed</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:color w:val="8496B0" w:themeColor="text2" w:themeTint="99"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="8496B0" w:themeColor="text2" w:themeTint="99"/>
        </w:rPr>
        <w:t>Comments page!!!</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:color w:val="8496B0" w:themeColor="text2" w:themeTint="99"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Scheduled work for next period of work</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>Andrew Campbell: New Comment form.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>David Doherty: Settings Page.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">John Lee Heaney: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Css</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>Matthew Mc Clean:  Add courses admin.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Ademoulu</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> - start topic page</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>** Additional Comments</w:t>
      </w:r>
      <w:bookmarkStart w:id="23" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="23"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)

Write-Output "Inserted Week6 section."
